$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '28.902.59'
$ws.Cells.Item(2, 5).Value = '  -2.33%  '
$ws.Cells.Item(3, 4).Value = '1.904.88'
$ws.Cells.Item(3, 5).Value = '  -4.53%  '
Set-TextValue 4 4 '1.006'
$ws.Cells.Item(4, 5).Value = '  +0.33%  '
Set-TextValue 5 4 '324.15'
$ws.Cells.Item(5, 5).Value = '  -0.80%  '
$ws.Cells.Item(6, 5).Value = '  +0.15%  '
$ws.Cells.Item(7, 5).Value = '  -2.23%  '
Set-TextValue 8 4 '0.3810'
$ws.Cells.Item(8, 5).Value = '  -3.70%  '
$ws.Cells.Item(9, 2).Value = 'OKB'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 9 4 '45.61'
$ws.Cells.Item(9, 5).Value = '  -2.05%  '
$ws.Cells.Item(10, 2).Value = 'Dogecoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 10 4 '0.07731'
$ws.Cells.Item(10, 5).Value = '  -2.91%  '
$ws.Cells.Item(11, 2).Value = 'Polygon'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 11 4 '0.9822'
$ws.Cells.Item(11, 5).Value = '  -2.11%  '
$ws.Cells.Item(12, 2).Value = 'Solana'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 12 4 '22.06'
$ws.Cells.Item(12, 5).Value = '  -3.91%  '
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.913.20'
$ws.Cells.Item(13, 5).Value = '  -3.61%  '
$ws.Cells.Item(14, 2).Value = 'Chainlink'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 14 4 '6.982'
$ws.Cells.Item(14, 5).Value = '  -4.12%  '
$ws.Cells.Item(15, 2).Value = 'Polkadot'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 15 4 '5.677'
$ws.Cells.Item(15, 5).Value = '  -3.67%  '
$ws.Cells.Item(16, 2).Value = 'TRON'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 16 4 '0.07049'
$ws.Cells.Item(16, 5).Value = '  -1.26%  '
$ws.Cells.Item(17, 2).Value = 'BinanceUSD'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 17 4 '1.005'
$ws.Cells.Item(17, 5).Value = '  +0.21%  '
$ws.Cells.Item(18, 2).Value = 'Litecoin'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 18 4 '84.13'
$ws.Cells.Item(18, 5).Value = '  -5.44%  '
$ws.Cells.Item(19, 2).Value = 'ShibaInu'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 19 4 '0.000009546'
$ws.Cells.Item(19, 5).Value = '  -4.54%  '
$ws.Cells.Item(20, 2).Value = 'Avalanche'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 20 4 '16.76'
$ws.Cells.Item(20, 5).Value = '  -4.06%  '
$ws.Cells.Item(21, 2).Value = 'Dai'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 21 4 '1.003'
$ws.Cells.Item(21, 5).Value = '  +0.28%  '
$ws.Cells.Item(22, 2).Value = 'WrappedBTC'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(22, 4).Value = '28.880.46'
$ws.Cells.Item(22, 5).Value = '  -2.45%  '
$ws.Cells.Item(23, 2).Value = 'Uniswap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 23 4 '5.332'
$ws.Cells.Item(23, 5).Value = '  -4.14%  '
$ws.Cells.Item(24, 2).Value = 'Cosmos'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 24 4 '10.93'
$ws.Cells.Item(24, 5).Value = '  -3.22%  '
$ws.Cells.Item(25, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(25, 4).Value = '2.139.23'
$ws.Cells.Item(25, 5).Value = '  -3.78%  '
$ws.Cells.Item(26, 2).Value = 'Toncoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 26 4 '2.082'
$ws.Cells.Item(26, 5).Value = '  -0.97%  '
$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 27 4 '156.77'
$ws.Cells.Item(27, 5).Value = '  -0.79%  '
$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 28 4 '19.17'
$ws.Cells.Item(28, 5).Value = '  -2.83%  '
$ws.Cells.Item(29, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 29 4 '5.590'
$ws.Cells.Item(29, 5).Value = '  -7.06%  '
$ws.Cells.Item(30, 2).Value = 'BitcoinCash'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 30 4 '117.81'
$ws.Cells.Item(30, 5).Value = '  -2.08%  '
$ws.Cells.Item(31, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 31 4 '1.841'
$ws.Cells.Item(31, 5).Value = '  -5.95%  '
$ws.Cells.Item(32, 2).Value = 'Stellar'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 32 4 '0.09269'
$ws.Cells.Item(32, 5).Value = '  -2.01%  '
$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 33 4 '0.8612'
$ws.Cells.Item(33, 5).Value = '  -5.84%  '
$ws.Cells.Item(34, 2).Value = 'Filecoin'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 34 4 '5.102'
$ws.Cells.Item(34, 5).Value = '  -3.24%  '
$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 35 4 '1.250'
$ws.Cells.Item(35, 5).Value = '  -7.74%  '
$ws.Cells.Item(36, 2).Value = 'HuobiToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 36 4 '3.015'
$ws.Cells.Item(36, 5).Value = '  -5.26%  '
$ws.Cells.Item(37, 2).Value = 'Hedera'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 37 4 '0.05710'
$ws.Cells.Item(37, 5).Value = '  -2.69%  '
$ws.Cells.Item(38, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 38 4 '1.145'
$ws.Cells.Item(38, 5).Value = '  -2.39%  '
$ws.Cells.Item(39, 2).Value = 'Frax'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 39 4 '1.004'
$ws.Cells.Item(39, 5).Value = '  +0.26%  '
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 40 4 '0.02035'
$ws.Cells.Item(40, 5).Value = '  -4.23%  '
$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 41 4 '7.488'
$ws.Cells.Item(41, 5).Value = '  -5.68%  '
$ws.Cells.Item(42, 2).Value = 'TheSandbox'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 42 4 '0.5513'
$ws.Cells.Item(42, 5).Value = '  -4.66%  '
$ws.Cells.Item(43, 2).Value = 'Algorand'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 43 4 '0.1757'
$ws.Cells.Item(43, 5).Value = '  -4.17%  '
$ws.Cells.Item(44, 2).Value = 'Aptos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 44 4 '9.325'
$ws.Cells.Item(44, 5).Value = '  -5.49%  '
$ws.Cells.Item(45, 2).Value = 'MXToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 45 4 '2.732'
$ws.Cells.Item(45, 5).Value = '  -0.96%  '
$ws.Cells.Item(46, 2).Value = 'Decentraland'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 46 4 '0.5201'
$ws.Cells.Item(46, 5).Value = '  -3.61%  '
$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 47 4 '11.28'
$ws.Cells.Item(47, 5).Value = '  -6.53%  '
$ws.Cells.Item(48, 2).Value = 'RenderToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 48 4 '2.089'
$ws.Cells.Item(48, 5).Value = '  -4.87%  '
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 49 4 '0.06815'
$ws.Cells.Item(49, 5).Value = '  -2.01%  '
Set-TextValue 50 4 '111.48'
$ws.Cells.Item(50, 5).Value = '  -2.59%  '
$ws.Cells.Item(51, 2).Value = 'PEPE'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 51 4 '0.000002568'
$ws.Cells.Item(51, 5).Value = '  -26.51%  '
